# Applies the "Deploying to gh-pages" metadata refresh to the
# StructureDefinition workbook:
#   - Metadata sheet: bump Version, Date, fill in Publisher, replace the
#     duplicated "Contact" rows with a single "Jurisdiction" row, and drop
#     the now-empty extra row.
#   - Elements sheet: give the root Extension row (row 2) its real
#     Short/Definition text instead of the generic placeholder.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# The "Contact" property was duplicated across rows 10 and 11; remove the
# duplicate so the remaining row can be repurposed for "Jurisdiction".
$meta.Rows.Item(11).Delete()

# Core metadata updates.
$meta.Cells.Item(3, 2).Value2 = "6.0.0"
$meta.Cells.Item(8, 2).Value2 = "2022-01-21T20:46:54+00:00"
$meta.Cells.Item(9, 2).Value2 = "Alvearie Team"

# Row 10 used to hold the leftover "Contact" property; repurpose it for
# "Jurisdiction".
$meta.Cells.Item(10, 1).Value2 = "Jurisdiction"
$meta.Cells.Item(10, 2).Value2 = "United States of America"

# Fill in the real Short/Definition text for the root Extension element.
$elements.Cells.Item(2, 11).Value2 = "Drug Coverage Indicator"
$elements.Cells.Item(2, 12).Value2 = "Indicates whether the member has drug benefit coverage: Y or N"
